$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.568.75'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.20%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.235.13'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '270.01'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '93.12'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +11.19%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.622'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.52%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.626'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.85%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '46.53'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +4.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0921'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.18'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +15.48%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.02%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.571.14'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.09'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.81%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.271.75'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.794'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.56%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.531.54'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.03%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.01'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.82%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '70.40'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.99%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.34'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.87%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '232.61'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.85'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -4.69%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +10.23%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.17'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +3.37%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.55'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +5.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '39.77'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.27'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '172.97'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0923'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.82'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.47'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.07%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.38%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -4.54%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -5.12%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.31'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -5.78%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.53'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +16.81%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -6.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.17'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.81%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.219'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +7.80%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '62.96'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.52%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.35'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.60%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0987'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.53%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '99.96'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -4.24%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.19'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.440'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.86%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.46'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -7.42%  '
